$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Rln3 -> Rxfp4, Target cluster: ECs)
$ws.Range("M2").Value = 7.889267333333334
$ws.Range("N2").Value = 23.667802
$ws.Range("O2").Value = 0.3739406795058886
$ws.Range("P2").Value = 0.3739406795058886
$ws.Range("Q2").Value = 5.563558659070667
$ws.Range("R2").Value = 50.072027931636
$ws.Range("S2").Value = 0.3739406795058886
$ws.Range("T2").Value = 0.3739406795058886

# Row 3 (Rln3 -> Rxfp4, Target cluster: FAPs)
$ws.Range("O3").Value = 0.1651501782221118
$ws.Range("P3").Value = 0.1651501782221118
$ws.Range("S3").Value = 0.1651501782221118
$ws.Range("T3").Value = 0.1651501782221118

# Row 4 (Rln3 -> Rxfp4, Target cluster: Inflammatory-Mac)
$ws.Range("M4").Value = 2.650898666666667
$ws.Range("N4").Value = 7.952696
$ws.Range("O4").Value = 0.1256490377156173
$ws.Range("P4").Value = 0.1256490377156173
$ws.Range("Q4").Value = 1.869429645125333
$ws.Range("R4").Value = 16.824866806128
$ws.Range("S4").Value = 0.1256490377156173
$ws.Range("T4").Value = 0.1256490377156173

# Row 5 (Rln3 -> Rxfp4, Target cluster: MuSCs)
$ws.Range("M5").Value = 4.010902
$ws.Range("N5").Value = 12.032706
$ws.Range("O5").Value = 0.1901113697813841
$ws.Range("P5").Value = 0.1901113697813841
$ws.Range("Q5").Value = 2.828512155812
$ws.Range("R5").Value = 25.456609402308
$ws.Range("S5").Value = 0.1901113697813841
$ws.Range("T5").Value = 0.1901113697813841

# Row 6 (Rln3 -> Rxfp4, Target cluster: Resolving-Mac)
$ws.Range("M6").Value = 3.062296333333334
$ws.Range("N6").Value = 9.186889000000001
$ws.Range("O6").Value = 0.1451487347749983
$ws.Range("P6").Value = 0.1451487347749983
$ws.Range("Q6").Value = 2.159549748044667
$ws.Range("R6").Value = 19.435947732402
$ws.Range("S6").Value = 0.1451487347749983
$ws.Range("T6").Value = 0.1451487347749983
